$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column A (shifts Features -> B, User Story -> C)
$ws.Range("A1").EntireColumn.Insert()

# New header for column A
$ws.Range("A1").Value = "Date"

# Existing data row (row 2) gets a date in the new column A
$ws.Range("A2").Value = 45319
$ws.Range("A2").NumberFormat = "mm-dd-yy"

# Insert two new rows (3 and 4) for the new backlog items - they inherit
# row 2's date formatting in column A automatically
$ws.Range("A3:A4").EntireRow.Insert()

$ws.Range("A3").Value = 45321
$ws.Range("B3").Value = "Data Model"
$ws.Range("C3").Value = "Update Data model so that each project has activities instead of only features. "

$ws.Range("A4").Value = 45321
$ws.Range("B4").Value = "Data Model"
$ws.Range("C4").Value = "Update Data model so that each user is part of a project, even if they have no currently assigned activity"

# Column widths (matches Excel's "best fit" auto-size after entering the data).
# Input values are tuned so the engine's rounded/stored width lands as close as
# possible to Excel's real best-fit pixel widths (9.453125 / 13.1796875 / 60.26953125 chars).
$ws.Range("A1").ColumnWidth = 8.6
$ws.Range("B1").ColumnWidth = 12.25
$ws.Range("C1").ColumnWidth = 59.42

# Move the selection to A5, matching where the cursor ends up after data entry
$ws.Range("A5").Select()
